# Append new job listing rows (4-10) to the tracking sheet, matching the
# target diff which extends the used range from A1:E3 to A1:E10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("2025-09-26 12:03", "Test Government Job - Software Engineer", "https://example.com/test-job", "TEST", "Not Applied"),
    @("2025-09-26 12:03", "⚫Advertisement No.37/2025⚫Advertisement No.36/2025⚫Advertisement No.35/2025", "https://www.ppsc.gop.pk/Adds/Advt No-37-2025 26-09-2025  X7 Version.pdf", "PPSC", "Not Applied"),
    @("2025-09-26 17:44", "Test Government Job - Software Engineer", "https://example.com/test-job", "TEST", "Not Applied"),
    @("2025-09-26 17:45", "⚫Advertisement No.37/2025⚫Advertisement No.36/2025⚫Advertisement No.35/2025", "https://www.ppsc.gop.pk/Adds/Advt No-37-2025 26-09-2025  X7 Version.pdf", "PPSC", "Not Applied"),
    @("2025-09-26 17:45", "Test Government Job - Software Engineer", "https://example.com/test-job", "TEST", "Not Applied"),
    @("2025-09-26 17:45", "⚫Advertisement No.37/2025⚫Advertisement No.36/2025⚫Advertisement No.35/2025", "https://www.ppsc.gop.pk/Adds/Advt No-37-2025 26-09-2025  X7 Version.pdf", "PPSC", "Not Applied"),
    @("2025-09-28 21:09", "⚫Advertisement No.37/2025⚫Advertisement No.36/2025⚫Advertisement No.35/2025", "https://www.ppsc.gop.pk/Adds/Advt No-37-2025 26-09-2025  X7 Version.pdf", "PPSC", "Not Applied")
)

$startRow = 4
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}

Write-Output "Added $($rows.Count) rows starting at row $startRow"
